$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13 (the lone "8452037 - Elisabeth..." row without a label) is removed
# entirely; all subsequent rows shift up by one.
$ws.Rows.Item(13).Delete()

# The teacher identification text now appears under "Objetivos:" (row 10),
# replacing the old objectives paragraph.
$teacher = "8452037 - Elisabeth Pinheiro da Silva Kondracki de Alcantara"
$ws.Range("B10").Value = $teacher
$ws.Range("C10").Value = $teacher

# "Programa resumido:" (now row 13 after the shift) becomes "Semestral".
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# "Programa:" (now row 15 after the shift) becomes the activation date
# text "01/01/2017" (the same literal text already used for "Ativação:" in
# B8/C8). Copy it over so it stays a plain text value instead of being
# reinterpreted as a date value.
$ws.Range("B8").Copy()
$ws.Range("B15").PasteSpecial(-4163)
$ws.Range("C8").Copy()
$ws.Range("C15").PasteSpecial(-4163)

# "Método:" (now row 18 after the shift) now shows the teacher identification.
$ws.Range("B18").Value = $teacher
$ws.Range("C18").Value = $teacher

# "Critério:" (now row 19 after the shift) gets the "A cada semestre..." text.
$metodoTexto = "A cada semestre é proposto um programa com cerca de 8 (oito) peças, sendo duas ou três de semestres anteriores e, consequentemente, cinco ou seis inéditas a ser apresentado pelo CORAL da EEL-USP em performances públicas definidas durante o período letivo."
$ws.Range("B19").Value = $metodoTexto
$ws.Range("C19").Value = $metodoTexto

# "Norma de recuperação:" (now row 20 after the shift) gets the "Sendo uma
# atividade prática..." text.
$criterioTexto = "Sendo uma atividade prática e de grupo, fica inviável a realização de provas ou outras formas similares de avaliação. Esta se dará no dia a dia do aluno, levando em conta: assiduidade, pontualidade e material completo na pasta; participação construtiva em sala de aula e nas apresentações públicas - prontidão, envolvimento e seu real aproveitamento vocal e musical."
$ws.Range("B20").Value = $criterioTexto
$ws.Range("C20").Value = $criterioTexto

# "Bibliografia:" (now row 21 after the shift) becomes "não tem".
$ws.Range("B21").Value = "não tem"
$ws.Range("C21").Value = "não tem"
